# Correzione nome foglio di lavoro, in TestCases
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TEST_CASE")
$ws.Name = "TestCases"

# Reflect the final active-cell selection recorded in the saved view
$ws.Range("D8").Select() | Out-Null
